# Apply the "LC" column insertion + TOST number-format tweak + sheet
# selection changes described by the commit.

$wb = $excel.ActiveWorkbook

$wsLcdod = $wb.Worksheets.Item("lcdod")
$wsHyst  = $wb.Worksheets.Item("hysteresis")

# --- lcdod: insert a new "LC" column right after "foam" -------------------
# Duplicate column C (the old "L135" column) into a freshly inserted column
# D; this shifts the old D/E columns to E/F and keeps the per-row styling
# (borders/fills) consistent with the rest of the table.
$wsLcdod.Columns.Item(3).Copy() | Out-Null
$wsLcdod.Columns.Item(4).Insert() | Out-Null

# The original "L135" header that is now duplicated in C1/D1 should read
# "LC" in the newly inserted (leftmost) slot, leaving D1 as "L135".
$wsLcdod.Range("C1").Value = "LC"

# --- styles: TOST values now print with 3 decimals instead of 2 ----------
$wsLcdod.Range("C2:F11").NumberFormat = "0.000"

# --- selections / active sheet --------------------------------------------
$wsHyst.Range("C19").Select() | Out-Null
$wsLcdod.Range("C25").Select() | Out-Null
$wsLcdod.Activate() | Out-Null
